$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump term version and publish date -------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.1.0
$ws1.Range("B3").Value = "1.1.0"

# Date: 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00
$ws1.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- Re-apply alignment/wrap so the "applyAlignment" flag is persisted ----
# (the header row + body rows already render with vertical=top/wrapText,
#  re-asserting WrapText makes the style record carry applyAlignment="true")
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A2:B14").WrapText = $true

$ws2 = $wb.Worksheets.Item("Include from CareSocialCodes")
$ws2.Range("A1:C1").WrapText = $true
$ws2.Range("A2:C2").WrapText = $true
$ws2.Range("A3:B4").WrapText = $true
